$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(38,8).Value = 257.4  # H38: 1335.1428 -> 257.4
$ws.Cells.Item(38,9).Value = 257.4  # I38: 472.83334 -> 257.4
$ws.Cells.Item(38,10).Value = 0  # J38: 6509 -> 0
$ws.Cells.Item(38,11).Value = 772.1999999999999  # K38: 1418.50002 -> 772.1999999999999
$ws.Cells.Item(38,12).Value = 0  # L38: 19527 -> 0
$ws.Cells.Item(38,13).Value = -400.1999999999999  # M38: -1046.50002 -> -400.1999999999999
$ws.Cells.Item(38,14).ClearContents()  # N38: was -20271
$ws.Cells.Item(41,8).Value = 561.5  # H41: 433.33334 -> 561.5
$ws.Cells.Item(41,9).Value = 625.5  # I41: 433.33334 -> 625.5
$ws.Cells.Item(41,10).Value = 497.5  # J41: 0 -> 497.5
$ws.Cells.Item(41,11).Value = 625.5  # K41: 433.33334 -> 625.5
$ws.Cells.Item(41,12).Value = 497.5  # L41: 0 -> 497.5
$ws.Cells.Item(41,13).Value = -185.5  # M41: 6.666659999999979 -> -185.5
$ws.Cells.Item(41,14).Value = -1377.5  # N41: <MISSING> -> -1377.5
$ws.Cells.Item(42,8).Value = 576.6667  # H42: 577.2222 -> 576.6667
$ws.Cells.Item(42,10).Value = 715  # J42: 715.8333 -> 715
$ws.Cells.Item(42,12).Value = 2145  # L42: 2147.4999 -> 2145
$ws.Cells.Item(42,14).Value = -2605  # N42: -2607.4999 -> -2605
$ws.Cells.Item(53,8).Value = 555  # H53: 517.0833 -> 555
$ws.Cells.Item(53,9).Value = 427.7143  # I53: 386.75 -> 427.7143
$ws.Cells.Item(53,11).Value = 427.7143  # K53: 386.75 -> 427.7143
$ws.Cells.Item(53,13).Value = 209.2857  # M53: 250.25 -> 209.2857
$ws.Cells.Item(58,8).Value = 8829.571  # H58: 8386.166999999999 -> 8829.571
$ws.Cells.Item(58,10).Value = 12181.4  # J58: 12354.25 -> 12181.4
$ws.Cells.Item(58,12).Value = 36544.2  # L58: 37062.75 -> 36544.2
$ws.Cells.Item(58,14).Value = -36844.2  # N58: -37362.75 -> -36844.2
$ws.Cells.Item(74,8).Value = 3000  # H74: 3099.6667 -> 3000
$ws.Cells.Item(74,9).Value = 3000  # I74: 2899.5 -> 3000
$ws.Cells.Item(74,10).Value = 0  # J74: 3500 -> 0
$ws.Cells.Item(74,11).Value = 3000  # K74: 2899.5 -> 3000
$ws.Cells.Item(74,12).Value = 0  # L74: 3500 -> 0
$ws.Cells.Item(74,13).Value = -2064  # M74: -1963.5 -> -2064
$ws.Cells.Item(74,14).ClearContents()  # N74: was -5372
$ws.Cells.Item(77,8).Value = 3000  # H77: 3099.6667 -> 3000
$ws.Cells.Item(77,9).Value = 3000  # I77: 2899.5 -> 3000
$ws.Cells.Item(77,10).Value = 0  # J77: 3500 -> 0
$ws.Cells.Item(77,11).Value = 15000  # K77: 14497.5 -> 15000
$ws.Cells.Item(77,12).Value = 0  # L77: 17500 -> 0
$ws.Cells.Item(77,13).Value = -10320  # M77: -9817.5 -> -10320
$ws.Cells.Item(77,14).ClearContents()  # N77: was -26860
$ws.Cells.Item(86,8).Value = 4156  # H86: 4188.4 -> 4156
$ws.Cells.Item(86,9).Value = 4246  # I86: 4495 -> 4246
$ws.Cells.Item(86,10).Value = 4111  # J86: 4111.75 -> 4111
$ws.Cells.Item(86,11).Value = 4246  # K86: 4495 -> 4246
$ws.Cells.Item(86,12).Value = 4111  # L86: 4111.75 -> 4111
$ws.Cells.Item(86,13).Value = -3123  # M86: -3372 -> -3123
$ws.Cells.Item(86,14).Value = -6357  # N86: -6357.75 -> -6357
$ws.Cells.Item(89,8).Value = 4156  # H89: 4188.4 -> 4156
$ws.Cells.Item(89,9).Value = 4246  # I89: 4495 -> 4246
$ws.Cells.Item(89,10).Value = 4111  # J89: 4111.75 -> 4111
$ws.Cells.Item(89,11).Value = 21230  # K89: 22475 -> 21230
$ws.Cells.Item(89,12).Value = 20555  # L89: 20558.75 -> 20555
$ws.Cells.Item(89,13).Value = -15614  # M89: -16859 -> -15614
$ws.Cells.Item(89,14).Value = -31787  # N89: -31790.75 -> -31787
$ws.Cells.Item(98,8).Value = 884  # H98: 942.5454999999999 -> 884
$ws.Cells.Item(98,9).Value = 860.8  # I98: 929.7778 -> 860.8
$ws.Cells.Item(98,11).Value = 860.8  # K98: 929.7778 -> 860.8
$ws.Cells.Item(98,13).Value = 637.2  # M98: 568.2222 -> 637.2
$ws.Cells.Item(122,8).Value = 884  # H122: 942.5454999999999 -> 884
$ws.Cells.Item(122,9).Value = 860.8  # I122: 929.7778 -> 860.8
$ws.Cells.Item(122,11).Value = 2582.4  # K122: 2789.3334 -> 2582.4
$ws.Cells.Item(122,13).Value = -132.3999999999996  # M122: -339.3334 -> -132.3999999999996
$ws.Cells.Item(138,8).Value = 3866.641  # H138: 3880.6843 -> 3866.641
$ws.Cells.Item(138,10).Value = 4162.6875  # J138: 4189.4517 -> 4162.6875
$ws.Cells.Item(138,12).Value = 12488.0625  # L138: 12568.3551 -> 12488.0625
$ws.Cells.Item(138,14).Value = -22768.0625  # N138: -22848.3551 -> -22768.0625

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32,8).Value = 5463.508  # H32: 5661.55 -> 5463.508
$ws.Cells.Item(32,9).Value = 5309.8687  # I32: 5418.5254 -> 5309.8687
$ws.Cells.Item(32,10).Value = 10149.5  # J32: 20000 -> 10149.5
$ws.Cells.Item(32,11).Value = 5309.8687  # K32: 5418.5254 -> 5309.8687
$ws.Cells.Item(32,12).Value = 10149.5  # L32: 20000 -> 10149.5
$ws.Cells.Item(32,13).Value = -5022.8687  # M32: -5131.5254 -> -5022.8687
$ws.Cells.Item(32,14).Value = -10723.5  # N32: -20574 -> -10723.5
$ws.Cells.Item(55,8).Value = 49026.5  # H55: 47221.2 -> 49026.5
$ws.Cells.Item(55,10).Value = 49026.5  # J55: 47221.2 -> 49026.5
$ws.Cells.Item(55,12).Value = 49026.5  # L55: 47221.2 -> 49026.5
$ws.Cells.Item(55,14).Value = -49656.5  # N55: -47851.2 -> -49656.5
$ws.Cells.Item(61,8).Value = 1999.75  # H61: 1658.3334 -> 1999.75
$ws.Cells.Item(61,9).Value = 1238  # I61: 979.46155 -> 1238
$ws.Cells.Item(61,11).Value = 1238  # K61: 979.46155 -> 1238
$ws.Cells.Item(61,13).Value = -1026  # M61: -767.46155 -> -1026
$ws.Cells.Item(136,8).Value = 1999.75  # H136: 1658.3334 -> 1999.75
$ws.Cells.Item(136,9).Value = 1238  # I136: 979.46155 -> 1238
$ws.Cells.Item(136,11).Value = 3714  # K136: 2938.38465 -> 3714
$ws.Cells.Item(136,13).Value = -1164  # M136: -388.38465 -> -1164

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(44,8).Value = 25000  # H44: 0 -> 25000
$ws.Cells.Item(44,10).Value = 25000  # J44: 0 -> 25000
$ws.Cells.Item(44,12).Value = 25000  # L44: 0 -> 25000
$ws.Cells.Item(44,14).Value = -25994  # N44: <MISSING> -> -25994
$ws.Cells.Item(80,8).Value = 834.0714  # H80: 792.625 -> 834.0714
$ws.Cells.Item(80,9).Value = 1309.5  # I80: 1051.6 -> 1309.5
$ws.Cells.Item(80,10).Value = 643.9  # J80: 674.9091 -> 643.9
$ws.Cells.Item(80,11).Value = 1309.5  # K80: 1051.6 -> 1309.5
$ws.Cells.Item(80,12).Value = 643.9  # L80: 674.9091 -> 643.9
$ws.Cells.Item(80,13).Value = -311.5  # M80: -53.59999999999991 -> -311.5
$ws.Cells.Item(80,14).Value = -2639.9  # N80: -2670.9091 -> -2639.9
$ws.Cells.Item(82,8).Value = 70283  # H82: 41761 -> 70283
$ws.Cells.Item(82,9).Value = 0  # I82: 27500 -> 0
$ws.Cells.Item(82,11).Value = 0  # K82: 27500 -> 0
$ws.Cells.Item(82,13).ClearContents()  # M82: was -27117
$ws.Cells.Item(83,8).Value = 834.0714  # H83: 792.625 -> 834.0714
$ws.Cells.Item(83,9).Value = 1309.5  # I83: 1051.6 -> 1309.5
$ws.Cells.Item(83,10).Value = 643.9  # J83: 674.9091 -> 643.9
$ws.Cells.Item(83,11).Value = 6547.5  # K83: 5258 -> 6547.5
$ws.Cells.Item(83,12).Value = 3219.5  # L83: 3374.5455 -> 3219.5
$ws.Cells.Item(83,13).Value = -1555.5  # M83: -266 -> -1555.5
$ws.Cells.Item(83,14).Value = -13203.5  # N83: -13358.5455 -> -13203.5
$ws.Cells.Item(85,8).Value = 70283  # H85: 41761 -> 70283
$ws.Cells.Item(85,9).Value = 0  # I85: 27500 -> 0
$ws.Cells.Item(85,11).Value = 0  # K85: 27500 -> 0
$ws.Cells.Item(85,13).ClearContents()  # M85: was -26174
$ws.Cells.Item(86,8).Value = 2691.3  # H86: 2617.4546 -> 2691.3
$ws.Cells.Item(86,9).Value = 2691.3  # I86: 2617.4546 -> 2691.3
$ws.Cells.Item(86,11).Value = 2691.3  # K86: 2617.4546 -> 2691.3
$ws.Cells.Item(86,13).Value = -1568.3  # M86: -1494.4546 -> -1568.3
$ws.Cells.Item(89,8).Value = 2691.3  # H89: 2617.4546 -> 2691.3
$ws.Cells.Item(89,9).Value = 2691.3  # I89: 2617.4546 -> 2691.3
$ws.Cells.Item(89,11).Value = 13456.5  # K89: 13087.273 -> 13456.5
$ws.Cells.Item(89,13).Value = -7840.5  # M89: -7471.273000000001 -> -7840.5
$ws.Cells.Item(94,8).Value = 1030.25  # H94: 1003.3333 -> 1030.25
$ws.Cells.Item(94,9).Value = 1030.25  # I94: 1003.3333 -> 1030.25
$ws.Cells.Item(94,11).Value = 1030.25  # K94: 1003.3333 -> 1030.25
$ws.Cells.Item(94,13).Value = -579.25  # M94: -552.3333 -> -579.25

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132,8).Value = 4165.5  # H132: 3057.0833 -> 4165.5
$ws.Cells.Item(132,9).Value = 3998.3333  # I132: 2631.889 -> 3998.3333
$ws.Cells.Item(132,11).Value = 11994.9999  # K132: 7895.667 -> 11994.9999
$ws.Cells.Item(132,13).Value = -9464.999899999999  # M132: -5365.667 -> -9464.999899999999
$ws.Cells.Item(141,8).Value = 55111  # H141: 58240.832 -> 55111
$ws.Cells.Item(141,10).Value = 55111  # J141: 58240.832 -> 55111
$ws.Cells.Item(141,12).Value = 55111  # L141: 58240.832 -> 55111
$ws.Cells.Item(141,14).Value = -65471  # N141: -68600.83199999999 -> -65471

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(3,8).Value = 0  # H3: 5000 -> 0
$ws.Cells.Item(3,9).Value = 0  # I3: 5000 -> 0
$ws.Cells.Item(3,11).Value = 0  # K3: 15000 -> 0
$ws.Cells.Item(3,13).ClearContents()  # M3: was -14888
$ws.Cells.Item(29,8).Value = 5037.5  # H29: 6732.3335 -> 5037.5
$ws.Cells.Item(29,9).Value = 0  # I29: 199 -> 0
$ws.Cells.Item(29,10).Value = 5037.5  # J29: 9999 -> 5037.5
$ws.Cells.Item(29,11).Value = 0  # K29: 597 -> 0
$ws.Cells.Item(29,12).Value = 15112.5  # L29: 29997 -> 15112.5
$ws.Cells.Item(29,13).ClearContents()  # M29: was -320
$ws.Cells.Item(29,14).Value = -15666.5  # N29: -30551 -> -15666.5
$ws.Cells.Item(34,8).Value = 4209.778  # H34: 4654.222 -> 4209.778
$ws.Cells.Item(34,10).Value = 5631.6665  # J34: 6298.3335 -> 5631.6665
$ws.Cells.Item(34,12).Value = 16894.9995  # L34: 18895.0005 -> 16894.9995
$ws.Cells.Item(34,14).Value = -17062.9995  # N34: -19063.0005 -> -17062.9995
$ws.Cells.Item(107,8).Value = 283.8  # H107: 278.16666 -> 283.8
$ws.Cells.Item(107,10).Value = 283.8  # J107: 278.16666 -> 283.8
$ws.Cells.Item(107,12).Value = 851.4000000000001  # L107: 834.4999799999999 -> 851.4000000000001
$ws.Cells.Item(107,14).Value = -4691.4  # N107: -4674.49998 -> -4691.4
$ws.Cells.Item(113,8).Value = 2889  # H113: 2958.2307 -> 2889
$ws.Cells.Item(113,10).Value = 2889  # J113: 2958.2307 -> 2889
$ws.Cells.Item(113,12).Value = 8667  # L113: 8874.6921 -> 8667
$ws.Cells.Item(113,14).Value = -13007  # N113: -13214.6921 -> -13007

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113,8).Value = 2233.3333  # H113: 1850 -> 2233.3333
$ws.Cells.Item(113,9).Value = 2233.3333  # I113: 1850 -> 2233.3333
$ws.Cells.Item(113,11).Value = 2233.3333  # K113: 1850 -> 2233.3333
$ws.Cells.Item(113,13).Value = -63.33329999999978  # M113: 320 -> -63.33329999999978
$ws.Cells.Item(131,8).Value = 0  # H131: 20000 -> 0
$ws.Cells.Item(131,9).Value = 0  # I131: 20000 -> 0
$ws.Cells.Item(131,11).Value = 0  # K131: 20000 -> 0
$ws.Cells.Item(131,13).ClearContents()  # M131: was -14960

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(82,8).Value = 1738  # H82: 2474 -> 1738
$ws.Cells.Item(82,9).Value = 1738  # I82: 2474 -> 1738
$ws.Cells.Item(82,11).Value = 1738  # K82: 2474 -> 1738
$ws.Cells.Item(82,13).Value = -1377  # M82: -2113 -> -1377
$ws.Cells.Item(85,8).Value = 1738  # H85: 2474 -> 1738
$ws.Cells.Item(85,9).Value = 1738  # I85: 2474 -> 1738
$ws.Cells.Item(85,11).Value = 1738  # K85: 2474 -> 1738
$ws.Cells.Item(85,13).Value = -490  # M85: -1226 -> -490
$ws.Cells.Item(132,8).Value = 4787.8  # H132: 4875.8887 -> 4787.8
$ws.Cells.Item(132,9).Value = 4269.143  # I132: 4314.8335 -> 4269.143
$ws.Cells.Item(132,11).Value = 12807.429  # K132: 12944.5005 -> 12807.429
$ws.Cells.Item(132,13).Value = -10277.429  # M132: -10414.5005 -> -10277.429
$ws.Cells.Item(136,8).Value = 3069.1482  # H136: 3132.111 -> 3069.1482
$ws.Cells.Item(136,10).Value = 2520.1667  # J136: 2803.5 -> 2520.1667
$ws.Cells.Item(136,12).Value = 7560.500100000001  # L136: 8410.5 -> 7560.500100000001
$ws.Cells.Item(136,14).Value = -12660.5001  # N136: -13510.5 -> -12660.5001

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(64,8).Value = 0  # H64: 114 -> 0
$ws.Cells.Item(64,10).Value = 0  # J64: 114 -> 0
$ws.Cells.Item(64,12).Value = 0  # L64: 114 -> 0
$ws.Cells.Item(64,14).ClearContents()  # N64: was -610
$ws.Cells.Item(67,8).Value = 0  # H67: 114 -> 0
$ws.Cells.Item(67,10).Value = 0  # J67: 114 -> 0
$ws.Cells.Item(67,12).Value = 0  # L67: 114 -> 0
$ws.Cells.Item(67,14).ClearContents()  # N67: was -1830
$ws.Cells.Item(81,8).Value = 7841.1113  # H81: 7155.8 -> 7841.1113
$ws.Cells.Item(81,9).Value = 4715.2  # I81: 4094 -> 4715.2
$ws.Cells.Item(81,11).Value = 9430.4  # K81: 8188 -> 9430.4
$ws.Cells.Item(81,13).Value = -8369.4  # M81: -7127 -> -8369.4
$ws.Cells.Item(84,8).Value = 7841.1113  # H84: 7155.8 -> 7841.1113
$ws.Cells.Item(84,9).Value = 4715.2  # I84: 4094 -> 4715.2
$ws.Cells.Item(84,11).Value = 47152  # K84: 40940 -> 47152
$ws.Cells.Item(84,13).Value = -41848  # M84: -35636 -> -41848
$ws.Cells.Item(132,8).Value = 2832.9656  # H132: 2878.6072 -> 2832.9656
$ws.Cells.Item(132,9).Value = 2109.158  # I132: 2139.9443 -> 2109.158
$ws.Cells.Item(132,11).Value = 6327.474  # K132: 6419.8329 -> 6327.474
$ws.Cells.Item(132,13).Value = -3797.474  # M132: -3889.8329 -> -3797.474
